# Rename AHB-Diff header columns from "_old"/"_new" suffixes to the
# format-version-specific suffixes "_FV2304"/"_FV2310" respectively, wrap the
# sheet's data range in an Excel Table ("Table1"), and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Rewrite the header row (row 1, columns A:U) -----------------------
$headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304",
    "diff",
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2) Turn the used range into a real Excel Table named "Table1" --------
$dataRange = $ws.Range("A1:U57")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)  # xlSrcRange = 1, xlYes = 1
$table.Name = "Table1"

# --- 3) Freeze the header row (row 1) --------------------------------------
[void]($ws.Range("A2").Select())
[void]($excel.ActiveWindow.FreezePanes = $true)
